# Update cryptos list prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain TEXT even when it looks like a
# number (e.g. "1.00", "151.55") instead of being auto-converted to a numeric
# cell. We flip the cell to text format, assign the value, then restore the
# default "Normal" style so no visible/structural formatting change remains.
function Set-TextValue {
    param($Worksheet, $Address, $NewValue)
    $Worksheet.Range($Address).NumberFormat = "@"
    $Worksheet.Range($Address).Value = $NewValue
    $Worksheet.Range($Address).Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '72.390.94'
$ws.Range('E2').Value = '  +0.82%  '
# Row 3
$ws.Range('D3').Value = '3.970.52'
# Row 4
$ws.Range('E4').Value = '  +0.13%  '
# Row 5
Set-TextValue $ws 'D5' '584.59'
$ws.Range('E5').Value = '  +9.32%  '
# Row 6
Set-TextValue $ws 'D6' '151.55'
$ws.Range('E6').Value = '  +1.67%  '
# Row 7
Set-TextValue $ws 'D7' '0.674'
$ws.Range('E7').Value = '  -2.84%  '
# Row 8
Set-TextValue $ws 'D8' '1.00'
$ws.Range('E8').Value = '  +0.10%  '
# Row 9
Set-TextValue $ws 'D9' '0.744'
$ws.Range('E9').Value = '  +0.12%  '
# Row 10
Set-TextValue $ws 'D10' '0.166'
$ws.Range('E10').Value = '  -1.83%  '
# Row 11
Set-TextValue $ws 'D11' '53.13'
$ws.Range('E11').Value = '  +5.63%  '
# Row 12
Set-TextValue $ws 'D12' '0.0000316'
$ws.Range('E12').Value = '  -1.82%  '
# Row 13
Set-TextValue $ws 'D13' '10.72'
$ws.Range('E13').Value = '  +0.50%  '
# Row 14
$ws.Range('D14').Value = '4.608.24'
$ws.Range('E14').Value = '  -0.92%  '
# Row 15
$ws.Range('D15').Value = '3.981.66'
$ws.Range('E15').Value = '  -0.57%  '
# Row 16
$ws.Range('E16').Value = '  +8.60%  '
# Row 17
Set-TextValue $ws 'D17' '13.94'
$ws.Range('E17').Value = '  -0.32%  '
# Row 18
Set-TextValue $ws 'D18' '20.37'
$ws.Range('E18').Value = '  -0.49%  '
# Row 19
$ws.Range('E19').Value = '  -0.29%  '
# Row 20
$ws.Range('D20').Value = '72.392.50'
$ws.Range('E20').Value = '  +0.79%  '
# Row 21
Set-TextValue $ws 'D21' '427.23'
$ws.Range('E21').Value = '  +0.11%  '
# Row 22
Set-TextValue $ws 'D22' '4.68'
$ws.Range('E22').Value = '  +11.32%  '
# Row 23
Set-TextValue $ws 'D23' '95.32'
$ws.Range('E23').Value = '  -1.58%  '
# Row 24
$ws.Range('E24').Value = '  -1.21%  '
# Row 25
Set-TextValue $ws 'D25' '4.46'
$ws.Range('E25').Value = '  +21.25%  '
# Row 26
Set-TextValue $ws 'D26' '14.17'
$ws.Range('E26').Value = '  -0.42%  '
# Row 27
Set-TextValue $ws 'D27' '11.22'
$ws.Range('E27').Value = '  +0.71%  '
# Row 28
Set-TextValue $ws 'D28' '10.60'
$ws.Range('E28').Value = '  -0.73%  '
# Row 29
$ws.Range('E29').Value = '  +1.11%  '
# Row 30
Set-TextValue $ws 'D30' '36.14'
$ws.Range('E30').Value = '  -1.49%  '
# Row 31
Set-TextValue $ws 'D31' '7.76'
$ws.Range('E31').Value = '  +5.38%  '
# Row 32
Set-TextValue $ws 'D32' '49.76'
$ws.Range('E32').Value = '  +4.15%  '
# Row 33
Set-TextValue $ws 'D33' '13.40'
$ws.Range('E33').Value = '  +0.57%  '
# Row 34
Set-TextValue $ws 'D34' '0.130'
$ws.Range('E34').Value = '  -0.13%  '
# Row 35
Set-TextValue $ws 'D35' '679.28'
$ws.Range('E35').Value = '  +0.79%  '
# Row 36
Set-TextValue $ws 'D36' '68.33'
$ws.Range('E36').Value = '  +4.42%  '
# Row 37
Set-TextValue $ws 'D37' '0.435'
$ws.Range('E37').Value = '  -1.83%  '
# Row 38
$ws.Range('D38').Value = '0.0₃0849'
$ws.Range('E38').Value = '  +4.34%  '
# Row 39
Set-TextValue $ws 'D39' '3.36'
$ws.Range('E39').Value = '  +0.23%  '
# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D40' '0.998'
$ws.Range('E40').Value = '  -0.06%  '
# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D41' '0.145'
$ws.Range('E41').Value = '  -2.79%  '
# Row 42
Set-TextValue $ws 'D42' '10.97'
$ws.Range('E42').Value = '  +11.81%  '
# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D43' '1.00'
$ws.Range('E43').Value = '  -0.18%  '
# Row 44
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws 'D44' '3.29'
$ws.Range('E44').Value = '  -3.60%  '
# Row 45
Set-TextValue $ws 'D45' '0.0484'
$ws.Range('E45').Value = '  -0.27%  '
# Row 46
$ws.Range('E46').Value = '  +4.25%  '
# Row 47
Set-TextValue $ws 'D47' '0.147'
$ws.Range('E47').Value = '  -0.93%  '
# Row 48
$ws.Range('E48').Value = '  +0.47%  '
# Row 49
Set-TextValue $ws 'D49' '3.42'
$ws.Range('E49').Value = '  +5.71%  '
# Row 50
Set-TextValue $ws 'D50' '2.98'
$ws.Range('E50').Value = '  -0.49%  '
# Row 51
$ws.Range('E51').Value = '  +6.76%  '
